# OysterCatalyst_PercentCover_Combined.xlsx edit
# Adds a "TotalType" classification column (Shell/Mud) and a matching set
# of new "Mud" rows (PC + Image J) beneath the existing "Shell" data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the new "TotalType" column at K.
#    Inserting a fresh column before K shifts the existing K data (the
#    TotalShell formulas/values) over to L, keeping their formulas intact.
# ---------------------------------------------------------------------
$ws.Columns("K").Insert()

# Column widths: K (new, "TotalType") = 19 chars, L keeps its old width.
$ws.Columns("K").ColumnWidth = 18.166666666666664

# ---------------------------------------------------------------------
# 2. Header row: K1 = "TotalType" (bold/center/wrap like the other
#    headers), L1 keeps its original "TotalShell" header text/format.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 11).Value = "TotalType"
$ws.Cells.Item(1, 11).Font.Bold = $true
$ws.Cells.Item(1, 11).HorizontalAlignment = -4108
$ws.Cells.Item(1, 11).WrapText = $true

# ---------------------------------------------------------------------
# 3. Existing data rows (2-28) are all "Shell" total-type rows.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 11).Value = "Shell"
}

# ---------------------------------------------------------------------
# 4. New "Mud" rows for the PC section (rows 29-40) and the Image J
#    section (rows 41-55). Only columns A, B, K, L are populated, same
#    as the source data being mirrored.
# ---------------------------------------------------------------------

# -- PC / "Mud" rows (29-40) --------------------------------------------------
$pc_mud = @(
    @("Guana_North", 83),
    @("Guana_North", 86),
    @("Guana_North", 74),
    @("Guana_North", 89),
    @("Guana_Mid", 16),
    @("Guana_Mid", 4),
    @("Guana_Mid", 19),
    @("Guana_Mid", 15),
    @("Guana_South", 7),
    @("Guana_South", 7),
    @("Guana_South", 24),
    @("Guana_South", 32)
)

$r = 29
foreach ($row in $pc_mud) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "PC"
    $ws.Cells.Item($r, 11).Value = "Mud"
    $ws.Cells.Item($r, 12).Value = $row[1]
    $r++
}

# -- Image J / "Mud" rows (41-55) --------------------------------------------
$imagej_mud = @(
    @("Guana_North", 36),
    @("Guana_North", 84),
    @("Guana_North", 68),
    @("Guana_North", 64),
    @("Guana_North", 76),
    @("Guana_Mid", 0),
    @("Guana_Mid", 8),
    @("Guana_Mid", 4),
    @("Guana_Mid", 40),
    @("Guana_Mid", 40),
    @("Guana_South", 12),
    @("Guana_South", 32),
    @("Guana_South", 60),
    @("Guana_South", 0),
    @("Guana_South", 0)
)

$r = 41
foreach ($row in $imagej_mud) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "Image J"
    $ws.Cells.Item($r, 11).Value = "Mud"
    $ws.Cells.Item($r, 12).Value = $row[1]
    $r++
}

# The "Image J" A-column cells (rows 41-55) use the same small Arial 9
# style as the existing Image J rows (14-28). Mirror that formatting by
# copying it from row 14's A cell.
$ws.Cells.Item(14, 1).Copy()
$ws.Range($ws.Cells.Item(41, 1), $ws.Cells.Item(55, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Restore the on-sheet selection to match the authored file.
# ---------------------------------------------------------------------
$ws.Range("D30").Select()
